$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '27.245.19'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -1.71%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.821.99'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -1.67%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.005'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -1.49%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '314.55'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.69%  '
$ws.Range('E6').Value = '  -1.15%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4277'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -2.10%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3683'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -2.65%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07236'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -2.67%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.8623'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -2.32%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '21.02'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -2.39%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.819.26'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -1.96%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '6.677'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -1.65%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.07132'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.00%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '5.311'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -3.26%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '89.01'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +0.78%  '
$ws.Range('E17').Value = '  -1.33%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000008865'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -1.71%  '
$ws.Range('E19').Value = '  -1.24%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '15.07'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -2.56%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '27.279.95'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.63%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.145'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -2.34%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '10.88'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -2.30%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.058.16'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -1.65%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.004'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -1.66%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '153.52'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -2.08%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '18.34'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -1.69%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.127'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +6.88%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.239'
$ws.Range('D29').Style = "Normal"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '116.30'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -3.70%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.08903'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -1.49%  '
$ws.Range('E32').Value = '  -1.93%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.7603'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -1.29%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.461'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -2.15%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.807'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -7.10%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.004'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -1.42%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.115'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -2.07%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01976'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.07%  '
$ws.Range('E39').Value = '  -0.25%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.907'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +1.44%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '7.149'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +3.00%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.1687'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +0.64%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.5060'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -2.38%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '8.640'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -0.79%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '10.63'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -1.19%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.4804'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +1.62%  '
$ws.Range('B47').Value = 'Quant'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '106.79'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -2.87%  '
$ws.Range('E48').Value = '  -1.41%  '
$ws.Range('E49').Value = '  -1.14%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.663'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -2.97%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.812'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.78%  '
